$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BS in Information Systems -> BS in Biological Sciences, is -> bio)
$ws.Range("A2").Value = "BS in Biological Sciences"
$ws.Range("C2").Value = "bio"
$ws.Range("D2").Value = "bio_0"

# Row 3 (GenEd stays, is -> bio)
$ws.Range("C3").Value = "bio"
$ws.Range("D3").Value = "bio_1"

# Row 4 (BS in Computer Science -> BS in Information Systems, cs -> is)
$ws.Range("A4").Value = "BS in Information Systems"
$ws.Range("C4").Value = "is"
$ws.Range("D4").Value = "is_0"

# Row 5 (GenEd stays, cs -> is)
$ws.Range("C5").Value = "is"
$ws.Range("D5").Value = "is_1"

# Row 6 (BS in Business Administration -> BS in Computer Science, ba -> cs)
$ws.Range("A6").Value = "BS in Computer Science"
$ws.Range("C6").Value = "cs"
$ws.Range("D6").Value = "cs_0"

# Row 7 (EY2022 Qatar Business Administration... -> GenEd, ba -> cs)
$ws.Range("A7").Value = "GenEd"
$ws.Range("C7").Value = "cs"
$ws.Range("D7").Value = "cs_1"

# Row 8 (BS in Biological Sciences -> BS in Business Administration, bs -> ba)
$ws.Range("A8").Value = "BS in Business Administration"
$ws.Range("C8").Value = "ba"
$ws.Range("D8").Value = "ba_0"

# Row 9 (GenEd -> EY2022 Qatar Business Administration..., bs -> ba)
$ws.Range("A9").Value = "EY2022 Qatar Business Administration - University Core Requirements"
$ws.Range("C9").Value = "ba"
$ws.Range("D9").Value = "ba_1"
